$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = -0.34420358511100346
$ws.Range("B1").Value = 0.34347514946853153
$ws.Range("A2").Value = -0.26626456541364441
$ws.Range("B2").Value = 0.26387716285928775
$ws.Range("A3").Value = -0.16092843528110023
$ws.Range("B3").Value = 0.16030506218652008
$ws.Range("A4").Value = -0.14830506229860418
$ws.Range("B4").Value = 0.14775644175070113
$ws.Range("A5").Value = -0.14175644217255545
$ws.Range("B5").Value = 0.14067205274854278
$ws.Range("A6").Value = -0.092945929885447587
$ws.Range("B6").Value = 0.092838676571668888
$ws.Range("A7").Value = -0.072838677085645287
$ws.Range("B7").Value = 0.072589236525216094
$ws.Range("A8").Value = -0.052589237044035286
$ws.Range("B8").Value = 0.052386982435687734
$ws.Range("A9").Value = -0.046386982881745809
$ws.Range("B9").Value = 0.046219321578270467
$ws.Range("A10").Value = -0.040219322029173554
$ws.Range("B10").Value = 0.040198319839930718
$ws.Range("A11").Value = -0.035698320283074736
$ws.Range("B11").Value = 0.035657437677937054
$ws.Range("A12").Value = -0.029657438130672453
$ws.Range("B12").Value = 0.029511345777819553
$ws.Range("A13").Value = -0.039151698769972931
$ws.Range("B13").Value = 0.039085052535743792
$ws.Range("A14").Value = -0.027085053028626405
$ws.Range("B14").Value = 0.027052895758604478
$ws.Range("A15").Value = -0.021052896219425854
$ws.Range("B15").Value = 0.021027594905570268
$ws.Range("A16").Value = -0.015027595367813618
$ws.Range("B16").Value = 0.015004372281703038
$ws.Range("A17").Value = -0.0090043727457951306
$ws.Range("B17").Value = 0.008999999518649382
$ws.Range("A18").Value = -0.067210584476178781
$ws.Range("B18").Value = 0.067165732136619027
$ws.Range("A19").Value = -0.058165732553151717
$ws.Range("B19").Value = 0.057839739869824935
$ws.Range("A20").Value = -0.048839740294665646
$ws.Range("B20").Value = 0.04877795875919233
$ws.Range("A21").Value = -0.039777959185750333
$ws.Range("B21").Value = 0.03970658401268512
$ws.Range("A22").Value = -0.09394872800643661
$ws.Range("B22").Value = 0.093635342773897179
$ws.Range("A23").Value = -0.084635343205984093
$ws.Range("B23").Value = 0.084126982846638221
$ws.Range("A24").Value = -0.042126983469888302
$ws.Range("B24").Value = 0.041999999373315333
$ws.Range("A25").Value = -0.041726123297046058
$ws.Range("B25").Value = 0.041706005476417118
$ws.Range("A26").Value = -0.035706005911045224
$ws.Range("B26").Value = 0.035684280212006314
$ws.Range("A27").Value = -0.029684280647558126
$ws.Range("B27").Value = 0.029622154512672694
$ws.Range("A28").Value = -0.023622154951302932
$ws.Range("B28").Value = 0.02359265284014711
$ws.Range("A29").Value = -0.011592653314099977
$ws.Range("B29").Value = 0.011588014442839523
$ws.Range("A30").Value = -0.059433006206596861
$ws.Range("B30").Value = 0.059087572134683519
$ws.Range("A31").Value = -0.055342034664356277
$ws.Range("B31").Value = 0.055236689724528887
$ws.Range("A32").Value = -0.0060007182902745981
$ws.Range("B32").Value = 0.0059999995514372628

$ws.Columns.Item(2).ColumnWidth = 13.8333

